$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.028910173402931
$ws.Range("D2").Value = 1.038138879877032
$ws.Range("E2").Value = 1.028817105104139
$ws.Range("F2").Value = 1.04491299270009
$ws.Range("I2").Value = 1.031649153012914
$ws.Range("J2").Value = 1.034059873731365
$ws.Range("K2").Value = 1.040927765599148
$ws.Range("L2").Value = 1.031632820015957
$ws.Range("M2").Value = 1.047682707642708
$ws.Range("N2").Value = 1.015364879796459

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.029861912281781
$ws.Range("D3").Value = 1.039060141937868
$ws.Range("E3").Value = 1.029624621685975
$ws.Range("F3").Value = 1.045982851177255
$ws.Range("I3").Value = 1.031778346725071
$ws.Range("J3").Value = 1.034652409876265
$ws.Range("K3").Value = 1.041658685233339
$ws.Range("L3").Value = 1.032248353038171
$ws.Range("M3").Value = 1.048563226489965
$ws.Range("N3").Value = 1.015561069486552

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.030477618998801
$ws.Range("D4").Value = 1.039656466401835
$ws.Range("E4").Value = 1.030147424307103
$ws.Range("F4").Value = 1.046675649200385
$ws.Range("I4").Value = 1.031859689918602
$ws.Range("J4").Value = 1.035035090525951
$ws.Range("K4").Value = 1.04213121930382
$ws.Range("L4").Value = 1.032646280229683
$ws.Range("M4").Value = 1.049132912718913
$ws.Range("N4").Value = 1.01568775280963

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.030736429830419
$ws.Range("D5").Value = 1.039907209602056
$ws.Range("E5").Value = 1.030367277722646
$ws.Range("F5").Value = 1.046967027071237
$ws.Range("I5").Value = 1.03189334628578
$ws.Range("J5").Value = 1.035195793924522
$ws.Range("K5").Value = 1.042329771257815
$ws.Range("L5").Value = 1.032813480848956
$ws.Range("M5").Value = 1.04937239160569
$ws.Range("N5").Value = 1.015740946850848

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.030779883395647
$ws.Range("D6").Value = 1.039949313297573
$ws.Range("E6").Value = 1.030404195985761
$ws.Range("F6").Value = 1.047015958037714
$ws.Range("I6").Value = 1.031898965631727
$ws.Range("J6").Value = 1.035222766414228
$ws.Range("K6").Value = 1.042363103034891
$ws.Range("L6").Value = 1.032841549384506
$ws.Range("M6").Value = 1.049412600131087
$ws.Range("N6").Value = 1.015749874624496

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.030481077369467
$ws.Range("D7").Value = 1.03965981665585
$ws.Range("E7").Value = 1.030150361736612
$ws.Range("F7").Value = 1.046679542113363
$ws.Range("I7").Value = 1.031860141760645
$ws.Range("J7").Value = 1.035037238543448
$ws.Range("K7").Value = 1.042133872764162
$ws.Range("L7").Value = 1.032648514719318
$ws.Range("M7").Value = 1.049136112716189
$ws.Range("N7").Value = 1.015688463841408

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.029231844935385
$ws.Range("D8").Value = 1.038450181390496
$ws.Range("E8").Value = 1.029089949343292
$ws.Range("F8").Value = 1.045274446545542
$ws.Range("I8").Value = 1.031693280842871
$ws.Range("J8").Value = 1.034260274728755
$ws.Range("K8").Value = 1.041174869717002
$ws.Range("L8").Value = 1.031840917096196
$ws.Range("M8").Value = 1.047980296756139
$ws.Range("N8").Value = 1.015431237612965

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.027029556124566
$ws.Range("D9").Value = 1.036320272835548
$ws.Range("E9").Value = 1.027223602365327
$ws.Range("F9").Value = 1.042802568995058
$ws.Range("I9").Value = 1.03138202197568
$ws.Range("J9").Value = 1.032885610837269
$ws.Range("K9").Value = 1.039481806195927
$ws.Range("L9").Value = 1.030415078214094
$ws.Range("M9").Value = 1.045943114441034
$ws.Range("N9").Value = 1.014975959225688

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.025560740107344
$ws.Range("D10").Value = 1.03490148129766
$ws.Range("E10").Value = 1.025980938630209
$ws.Range("F10").Value = 1.041157440085454
$ws.Range("I10").Value = 1.03116297527142
$ws.Range("J10").Value = 1.03196547974663
$ws.Range("K10").Value = 1.038351009346728
$ws.Range("L10").Value = 1.029462720256051
$ws.Range("M10").Value = 1.044584712374935
$ws.Range("N10").Value = 1.014671103920854

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.024924585496776
$ws.Range("D11").Value = 1.034287412177596
$ws.Range("E11").Value = 1.025443236447465
$ws.Range("F11").Value = 1.040445751590834
$ws.Range("I11").Value = 1.03106539692451
$ws.Range("J11").Value = 1.031566187067923
$ws.Range("K11").Value = 1.037860876213792
$ws.Range("L11").Value = 1.029049922002814
$ws.Range("M11").Value = 1.04399645187341
$ws.Range("N11").Value = 1.014538784630145

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.024688267624654
$ws.Range("D12").Value = 1.034059361995853
$ws.Range("E12").Value = 1.025243567920642
$ws.Range("F12").Value = 1.040181499054182
$ws.Range("I12").Value = 1.031028742510847
$ws.Range("J12").Value = 1.031417741966989
$ws.Range("K12").Value = 1.03767874592897
$ws.Range("L12").Value = 1.028896527841701
$ws.Range("M12").Value = 1.043777936849329
$ws.Range("N12").Value = 1.014489588289573

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.0247389595864
$ws.Range("D13").Value = 1.034108277600832
$ws.Range("E13").Value = 1.025286394856581
$ws.Range("F13").Value = 1.040238177573488
$ws.Range("I13").Value = 1.031036623526373
$ws.Range("J13").Value = 1.031449589818238
$ws.Range("K13").Value = 1.03771781678635
$ws.Range("L13").Value = 1.028929434226932
$ws.Range("M13").Value = 1.043824809455302
$ws.Range("N13").Value = 1.0145001431974

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.02490505182836
$ws.Range("D14").Value = 1.03426856060815
$ws.Range("E14").Value = 1.025426730595699
$ws.Range("F14").Value = 1.040423906344901
$ws.Range("I14").Value = 1.031062375408702
$ws.Range("J14").Value = 1.031553919198166
$ws.Range("K14").Value = 1.037845822753137
$ws.Range("L14").Value = 1.029037243664557
$ws.Range("M14").Value = 1.043978389530797
$ws.Range("N14").Value = 1.014534719003948

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.025007383874895
$ws.Range("D15").Value = 1.034367321905811
$ws.Range("E15").Value = 1.025513203792303
$ws.Range("F15").Value = 1.04053835327451
$ws.Range("I15").Value = 1.031078187753204
$ws.Range("J15").Value = 1.03161818276302
$ws.Range("K15").Value = 1.037924681799369
$ws.Range("L15").Value = 1.029103660355958
$ws.Range("M15").Value = 1.044073014126539
$ws.Range("N15").Value = 1.014556016070455

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.025602956522667
$ws.Range("D16").Value = 1.034942240948384
$ws.Range("E16").Value = 1.02601663223277
$ws.Range("F16").Value = 1.041204686532692
$ws.Range("I16").Value = 1.031169393767946
$ws.Range("J16").Value = 1.031991961216801
$ws.Range("K16").Value = 1.038383527591035
$ws.Range("L16").Value = 1.02949010751021
$ws.Range("M16").Value = 1.044623752006005
$ws.Range("N16").Value = 1.014679878902581

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.02597650431579
$ws.Range("D17").Value = 1.035302947209648
$ws.Range("E17").Value = 1.026332521883443
$ws.Range("F17").Value = 1.041622837775095
$ws.Range("I17").Value = 1.031225874528443
$ws.Range("J17").Value = 1.032226190154131
$ws.Range("K17").Value = 1.038671218570722
$ws.Range("L17").Value = 1.029732403416175
$ws.Range("M17").Value = 1.044969198803245
$ws.Range("N17").Value = 1.014757490646246

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.026194374154366
$ws.Range("D18").Value = 1.035513367739179
$ws.Range("E18").Value = 1.026516811513554
$ws.Range("F18").Value = 1.041866802266019
$ws.Range("I18").Value = 1.03125855538421
$ws.Range("J18").Value = 1.032362727959061
$ws.Range("K18").Value = 1.0388389764664
$ws.Range("L18").Value = 1.029873689852075
$ws.Range("M18").Value = 1.045170685911073
$ws.Range("N18").Value = 1.014802729860228

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.02626865970504
$ws.Range("D19").Value = 1.035585120200683
$ws.Range("E19").Value = 1.026579655682371
$ws.Range("F19").Value = 1.041949998680082
$ws.Range("I19").Value = 1.031269654012127
$ws.Range("J19").Value = 1.032409269558804
$ws.Range("K19").Value = 1.038896169514883
$ws.Range("L19").Value = 1.029921857955249
$ws.Range("M19").Value = 1.045239386742461
$ws.Range("N19").Value = 1.01481815009857

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.02593642764539
$ws.Range("D20").Value = 1.035264244076895
$ws.Range("E20").Value = 1.026298626126521
$ws.Range("F20").Value = 1.041577967485301
$ws.Range("I20").Value = 1.031219841919122
$ws.Range("J20").Value = 1.032201068286559
$ws.Range("K20").Value = 1.038640356947457
$ws.Range("L20").Value = 1.029706411559488
$ws.Range("M20").Value = 1.044932136264001
$ws.Range("N20").Value = 1.014749166783551

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.024856142398992
$ws.Range("D21").Value = 1.034221360079613
$ws.Range("E21").Value = 1.025385403609074
$ws.Range("F21").Value = 1.040369211078298
$ws.Range("I21").Value = 1.031054803419576
$ws.Range("J21").Value = 1.031523200380066
$ws.Range("K21").Value = 1.037808130196429
$ws.Range("L21").Value = 1.029005498191962
$ws.Range("M21").Value = 1.043933164266433
$ws.Range("N21").Value = 1.014524538587911

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.024176798711855
$ws.Range("D22").Value = 1.033565903410124
$ws.Range("E22").Value = 1.024811560456231
$ws.Range("F22").Value = 1.039609798316963
$ws.Range("I22").Value = 1.030948668398208
$ws.Range("J22").Value = 1.031096245980742
$ws.Range("K22").Value = 1.037284454098021
$ws.Range("L22").Value = 1.02856444483429
$ws.Range("M22").Value = 1.04330501985387
$ws.Range("N22").Value = 1.014383033747334

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.024536943399756
$ws.Range("D23").Value = 1.033913349813454
$ws.Range("E23").Value = 1.025115733379108
$ws.Range("F23").Value = 1.040012322171776
$ws.Range("I23").Value = 1.031005156904187
$ws.Range("J23").Value = 1.031322653642219
$ws.Range("K23").Value = 1.037562104486593
$ws.Range("L23").Value = 1.028798289567707
$ws.Range("M23").Value = 1.043638015717279
$ws.Range("N23").Value = 1.014458073857128

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.025954536611684
$ws.Range("D24").Value = 1.035281732274227
$ws.Range("E24").Value = 1.026313942047194
$ws.Range("F24").Value = 1.041598242240216
$ws.Range("I24").Value = 1.031222568609541
$ws.Range("J24").Value = 1.032212420036613
$ws.Range("K24").Value = 1.038654302132994
$ws.Range("L24").Value = 1.02971815628639
$ws.Range("M24").Value = 1.044948883249167
$ws.Range("N24").Value = 1.014752928072751

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.027599013111158
$ws.Range("D25").Value = 1.036870707068546
$ws.Range("E25").Value = 1.027705826208965
$ws.Range("F25").Value = 1.043441120276615
$ws.Range("I25").Value = 1.031464527088701
$ws.Range("J25").Value = 1.033241647858485
$ws.Range("K25").Value = 1.03991987522174
$ws.Range("L25").Value = 1.030415078214094
$ws.Range("M25").Value = 1.04646982774578
$ws.Range("N25").Value = 1.015093896380358

